$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.788.26"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "2.083.92"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'233.41"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "'58.83"
$ws.Range("E7").Value = "  +3.08%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +1.95%  "
$ws.Range("D10").Value = "'0.0790"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("D12").Value = "2.389.91"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "'14.73"
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("D14").Value = "'21.19"
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D15").Value = "'0.776"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("D17").Value = "2.097.54"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "37.709.46"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'71.80"
$ws.Range("E20").Value = "  +1.71%  "
$ws.Range("D21").Value = "0.0₃0847"
$ws.Range("E21").Value = "  +3.31%  "
$ws.Range("D22").Value = "'228.25"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("E25").Value = "  +1.59%  "
$ws.Range("D26").Value = "'9.60"
$ws.Range("E26").Value = "  +7.81%  "
$ws.Range("D27").Value = "'171.27"
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("E28").Value = "  -2.22%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("E31").Value = "  +2.39%  "
$ws.Range("D32").Value = "'4.74"
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("E33").Value = "  +1.64%  "
$ws.Range("D34").Value = "'4.68"
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'3.43"
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").Value = "'17.31"
$ws.Range("E41").Value = "  +10.85%  "
$ws.Range("D42").Value = "'99.04"
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("D43").Value = "'0.0220"
$ws.Range("E43").Value = "  +2.77%  "
$ws.Range("D44").Value = "'2.92"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("D45").Value = "1.451.01"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("D47").Value = "'4.17"
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").Value = "2.276.01"
$ws.Range("E51").Value = "  +0.39%  "
